$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new column before column N (14th column), shifting N->O, O->P, P->Q
$ws.Columns("N").Insert()

# Move the active selection to match the target workbook state
$ws.Range("R8").Select()
